$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerated s_vals data (filtering save games) -- updated B:E inputs and
# recomputed G (sum) for each row.
$data = @{
    2 = @{ B = 0.6606524410359556;  C = 1.655778082260271;  D = 0.7527432677738641; E = 0.4942365360607697; G = 3.56341032713086 }
    3 = @{ B = 3.286832544864788;   C = 1.655778082260271;  D = 0.7527432677738641; E = 0.4942365360607697; G = 6.189590430959694 }
    4 = @{ B = 0.01293466051926884; C = 0.002571899574220771; D = 0.7527432677738641; E = 0.4942365360607697; G = 1.262486363928123 }
    5 = @{ B = 0.6606524410359556;  C = 1.655778082260271;  D = 0.1494219747398047; E = 0.4942365360607697; G = 2.960089034096801 }
    6 = @{ B = 0.6606524410359556;  C = 0.306821227259698;  D = 0.7527432677738641; E = 0.4942365360607697; G = 2.214453472130288 }
    7 = @{ B = 3.286832544864788;   C = 1.655778082260271;  D = 0.1494219747398047; E = 0.4942365360607697; G = 5.586269137925634 }
    8 = @{ B = 1.455362044514542;   C = 1.655778082260271;  D = 0.1494219747398047; E = 0.4942365360607697; G = 3.754798637575387 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
